# Apply the "commit code release v1" edit to ds-noidungkhac.xlsx
#
# Summary of the change:
#  - The date-range subtitle (row 4) is updated.
#  - The first data row (old row 7: STT=1, TH Hoa Binh / Lop 1, 11/07/2018,
#    "dsadwa"/"dwasdwada", subject list "Tin hoc, Ki thuat, Am nhac, Giao
#    trinh Tieu hoc") is removed entirely; the remaining data rows shift up
#    by one and the STT (sequence number) column is renumbered 1..5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the subtitle date range (merged cell A4:I4)
$ws.Range("A4").Value = "Từ ngày 22/07/2018 tới ngày 22/08/2018"

# 2. Remove the obsolete first data row (old row 7); rows below shift up
$ws.Rows("7").Delete()

# 3. Renumber the STT column for the now-shifted data rows (A7:A11 -> 1..5)
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
